$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added as row 33; existing rows 33-47
# shift down to 34-48 (dimension grows from A1:R47 to A1:R48).
$ws.Rows("33").Insert()

$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44553
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 95
$ws.Range("K33").Value = 45000
$ws.Range("L33").Value = 45000
$ws.Range("M33").Value = 45000
$ws.Range("N33").Value = "$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 1800
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
